$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row (2-369).
# All of them move forward by one day: 45180 -> 45181 (2023-09-11 -> 2023-09-12).
$ws.Range("C2:C369").Value = 45181
